$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for a new "No VA PMB" column at A, shifting the existing ----
# A:G data/header block over to B:H. Plain value-assignment instead of a
# real EntireColumn insert keeps the <cols> width metadata untouched (the
# original template already carries a stray gap at column 9/10).
for ($r = 1; $r -le 3; $r++) {
    for ($c = 7; $c -ge 1; $c--) {
        $ws.Cells.Item($r, $c + 1).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Carry the per-cell formatting along with it. A direct Copy+PasteSpecial
# straight onto the overlapping A1:G3 -> B1:H3 destination corrupts cells
# where source/target intersect, so stage through a scratch range well
# outside the used range first.
$ws.Range("A1:G3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1:AF3").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("Z1:AF3").Clear()

# --- New column A: "No VA PMB" header + the two VA numbers -----------------
$hdr = $ws.Range("A1")
$hdr.Value = "No VA PMB"
$hdr.Font.Bold = $true

# Leading apostrophe forces these long digit strings to be stored as text
# (quote-prefixed), matching how Excel keeps a 16-digit VA number from being
# mangled into a float.
$ws.Range("A2").Value = "'8257062100000000"
$ws.Range("A3").Value = "'8257062100000022"

# --- Selection left by the author after making the edit --------------------
$ws.Range("E13").Select()

Write-Host "done"
